$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.036019130622163
$ws.Cells.Item(2, 4).Value = 1.044334344382848
$ws.Cells.Item(2, 5).Value = 1.039666840014316
$ws.Cells.Item(2, 6).Value = 1.053352522305984
$ws.Cells.Item(2, 9).Value = 1.042660216335966
$ws.Cells.Item(2, 10).Value = 1.041129733014365
$ws.Cells.Item(2, 11).Value = 1.04710568636298
$ws.Cells.Item(2, 12).Value = 1.042451377914911
$ws.Cells.Item(2, 13).Value = 1.056098727274549
$ws.Cells.Item(2, 14).Value = 1.017656945390353

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.036847179549128
$ws.Cells.Item(3, 4).Value = 1.044980168692938
$ws.Cells.Item(3, 5).Value = 1.040442942468631
$ws.Cells.Item(3, 6).Value = 1.054135846360394
$ws.Cells.Item(3, 9).Value = 1.042883485650861
$ws.Cells.Item(3, 10).Value = 1.041602184318099
$ws.Cells.Item(3, 11).Value = 1.04756315935154
$ws.Cells.Item(3, 12).Value = 1.043037835863186
$ws.Cells.Item(3, 13).Value = 1.056695155072948
$ws.Cells.Item(3, 14).Value = 1.017814790389635

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.037383593322354
$ws.Cells.Item(4, 4).Value = 1.045398559577856
$ws.Cells.Item(4, 5).Value = 1.040946073719258
$ws.Cells.Item(4, 6).Value = 1.054643552688489
$ws.Cells.Item(4, 9).Value = 1.043026958912061
$ws.Cells.Item(4, 10).Value = 1.04190783271057
$ws.Cells.Item(4, 11).Value = 1.047858986349689
$ws.Cells.Item(4, 12).Value = 1.043417594318963
$ws.Cells.Item(4, 13).Value = 1.057081266371403
$ws.Cells.Item(4, 14).Value = 1.017916872298475

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.037609246220832
$ws.Cells.Item(5, 4).Value = 1.04557456886101
$ws.Cells.Item(5, 5).Value = 1.041157813510983
$ws.Cells.Item(5, 6).Value = 1.054857192345691
$ws.Cells.Item(5, 9).Value = 1.043087035572354
$ws.Cells.Item(5, 10).Value = 1.042036311853655
$ws.Cells.Item(5, 11).Value = 1.047983305543082
$ws.Cells.Item(5, 12).Value = 1.043577310522148
$ws.Cells.Item(5, 13).Value = 1.057243629674488
$ws.Cells.Item(5, 14).Value = 1.017959774051485

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.037647142763406
$ws.Cells.Item(6, 4).Value = 1.045604128447105
$ws.Cells.Item(6, 5).Value = 1.041193378600164
$ws.Cells.Item(6, 6).Value = 1.054893075050194
$ws.Cells.Item(6, 9).Value = 1.043097108639395
$ws.Cells.Item(6, 10).Value = 1.042057883124713
$ws.Cells.Item(6, 11).Value = 1.048004176518841
$ws.Cells.Item(6, 12).Value = 1.043604131387025
$ws.Cells.Item(6, 13).Value = 1.057270893610258
$ws.Cells.Item(6, 14).Value = 1.017966976640436

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.037386607941534
$ws.Cells.Item(7, 4).Value = 1.045400910962073
$ws.Cells.Item(7, 5).Value = 1.040948902121215
$ws.Cells.Item(7, 6).Value = 1.054646406568981
$ws.Cells.Item(7, 9).Value = 1.043027762601445
$ws.Cells.Item(7, 10).Value = 1.041909549517434
$ws.Cells.Item(7, 11).Value = 1.047860647693116
$ws.Cells.Item(7, 12).Value = 1.043419728198238
$ws.Cells.Item(7, 13).Value = 1.057083435714351
$ws.Cells.Item(7, 14).Value = 1.017917445607358

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.03629884617225
$ws.Cells.Item(8, 4).Value = 1.044552499255253
$ws.Cells.Item(8, 5).Value = 1.03992893148888
$ws.Cells.Item(8, 6).Value = 1.053617074547304
$ws.Cells.Item(8, 9).Value = 1.042735877161588
$ws.Cells.Item(8, 10).Value = 1.041289411637084
$ws.Cells.Item(8, 11).Value = 1.047260329760981
$ws.Cells.Item(8, 12).Value = 1.042649515043749
$ws.Cells.Item(8, 13).Value = 1.05630025405933
$ws.Cells.Item(8, 14).Value = 1.017710300821635

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.03438682229529
$ws.Cells.Item(9, 4).Value = 1.043061393350492
$ws.Cells.Item(9, 5).Value = 1.038138903449267
$ws.Cells.Item(9, 6).Value = 1.051809805552198
$ws.Cells.Item(9, 9).Value = 1.042213939518331
$ws.Cells.Item(9, 10).Value = 1.040196252276669
$ws.Cells.Item(9, 11).Value = 1.046201106723335
$ws.Cells.Item(9, 12).Value = 1.041294519599645
$ws.Cells.Item(9, 13).Value = 1.054921657791404
$ws.Cells.Item(9, 14).Value = 1.017344889564981

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.033115433498525
$ws.Cells.Item(10, 4).Value = 1.042070054615201
$ws.Cells.Item(10, 5).Value = 1.036950565292859
$ws.Cells.Item(10, 6).Value = 1.05060947780251
$ws.Cells.Item(10, 9).Value = 1.041860922321768
$ws.Cells.Item(10, 10).Value = 1.039467288044917
$ws.Cells.Item(10, 11).Value = 1.045494104871556
$ws.Cells.Item(10, 12).Value = 1.040392771043863
$ws.Cells.Item(10, 13).Value = 1.054003677842676
$ws.Cells.Item(10, 14).Value = 1.017101043889363

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.03256571023262
$ws.Cells.Item(11, 4).Value = 1.041641464598297
$ws.Cells.Item(11, 5).Value = 1.036437214666975
$ws.Cells.Item(11, 6).Value = 1.050090819112125
$ws.Cells.Item(11, 9).Value = 1.041706872596261
$ws.Cells.Item(11, 10).Value = 1.039151609070318
$ws.Cells.Item(11, 11).Value = 1.045187779530772
$ws.Cells.Item(11, 12).Value = 1.040002698151171
$ws.Cells.Item(11, 13).Value = 1.053606459816673
$ws.Cells.Item(11, 14).Value = 1.016995405406749

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.032361639871455
$ws.Cells.Item(12, 4).Value = 1.041482369105863
$ws.Cells.Item(12, 5).Value = 1.036246716742797
$ws.Cells.Item(12, 6).Value = 1.049898332124748
$ws.Cells.Item(12, 9).Value = 1.041649473571261
$ws.Cells.Item(12, 10).Value = 1.039034348220595
$ws.Cells.Item(12, 11).Value = 1.045073969619665
$ws.Cells.Item(12, 12).Value = 1.039857867865653
$ws.Cells.Item(12, 13).Value = 1.053458958057676
$ws.Cells.Item(12, 14).Value = 1.016956159276743

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.032405408139013
$ws.Cells.Item(13, 4).Value = 1.041516490987408
$ws.Cells.Item(13, 5).Value = 1.03628757086615
$ws.Cells.Item(13, 6).Value = 1.049939613695385
$ws.Cells.Item(13, 9).Value = 1.041661793902949
$ws.Cells.Item(13, 10).Value = 1.039059501214265
$ws.Cells.Item(13, 11).Value = 1.045098383434129
$ws.Cells.Item(13, 12).Value = 1.039888931706773
$ws.Cells.Item(13, 13).Value = 1.0534905957271
$ws.Cells.Item(13, 14).Value = 1.016964578027337

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.032548839227985
$ws.Cells.Item(14, 4).Value = 1.041628311625469
$ws.Cells.Item(14, 5).Value = 1.036421464297638
$ws.Cells.Item(14, 6).Value = 1.050074904680878
$ws.Cells.Item(14, 9).Value = 1.041702131603292
$ws.Cells.Item(14, 10).Value = 1.039141916318389
$ws.Cells.Item(14, 11).Value = 1.045178372511908
$ws.Cells.Item(14, 12).Value = 1.039990725198577
$ws.Cells.Item(14, 13).Value = 1.053594266394703
$ws.Cells.Item(14, 14).Value = 1.016992161455463

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.03263722798411
$ws.Cells.Item(15, 4).Value = 1.041697221571708
$ws.Cells.Item(15, 5).Value = 1.036503984813614
$ws.Cells.Item(15, 6).Value = 1.050158283971464
$ws.Cells.Item(15, 9).Value = 1.04172696141691
$ws.Cells.Item(15, 10).Value = 1.039192694542372
$ws.Cells.Item(15, 11).Value = 1.045227652879534
$ws.Cells.Item(15, 12).Value = 1.040053451545472
$ws.Cells.Item(15, 13).Value = 1.05365814702188
$ws.Cells.Item(15, 14).Value = 1.017009155563506

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.033151933744906
$ws.Cells.Item(16, 4).Value = 1.042098512915733
$ws.Cells.Item(16, 5).Value = 1.036984660290786
$ws.Cells.Item(16, 6).Value = 1.050643922632529
$ws.Cells.Item(16, 9).Value = 1.041871121073469
$ws.Cells.Item(16, 10).Value = 1.039488238040344
$ws.Cells.Item(16, 11).Value = 1.045514430810242
$ws.Cells.Item(16, 12).Value = 1.040418667248753
$ws.Cells.Item(16, 13).Value = 1.054030045793537
$ws.Cells.Item(16, 14).Value = 1.017108053718728

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.033475009490369
$ws.Cells.Item(17, 4).Value = 1.042350412029406
$ws.Cells.Item(17, 5).Value = 1.037286499993909
$ws.Cells.Item(17, 6).Value = 1.050948844485201
$ws.Cells.Item(17, 9).Value = 1.04196123035089
$ws.Cells.Item(17, 10).Value = 1.039673616832081
$ws.Cells.Item(17, 11).Value = 1.045694269405828
$ws.Cells.Item(17, 12).Value = 1.040647862956755
$ws.Cells.Item(17, 13).Value = 1.054263402404925
$ws.Cells.Item(17, 14).Value = 1.017170076412003

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.033663530842288
$ws.Cells.Item(18, 4).Value = 1.042497404676129
$ws.Cells.Item(18, 5).Value = 1.037462674376761
$ws.Cells.Item(18, 6).Value = 1.051126805387989
$ws.Cells.Item(18, 9).Value = 1.042013674602818
$ws.Cells.Item(18, 10).Value = 1.039781741908499
$ws.Cells.Item(18, 11).Value = 1.045799147902898
$ws.Cells.Item(18, 12).Value = 1.040781586479643
$ws.Cells.Item(18, 13).Value = 1.054399541632969
$ws.Cells.Item(18, 14).Value = 1.017206248180061

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.033727824716231
$ws.Cells.Item(19, 4).Value = 1.042547536196081
$ws.Cells.Item(19, 5).Value = 1.037522764971978
$ws.Cells.Item(19, 6).Value = 1.051187503230127
$ws.Cells.Item(19, 9).Value = 1.042031537212182
$ws.Cells.Item(19, 10).Value = 1.039818609177377
$ws.Cells.Item(19, 11).Value = 1.045834905613235
$ws.Cells.Item(19, 12).Value = 1.040827189066666
$ws.Cells.Item(19, 13).Value = 1.054445965996798
$ws.Cells.Item(19, 14).Value = 1.01721858095241

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.033440338571229
$ws.Cells.Item(20, 4).Value = 1.04232337898932
$ws.Cells.Item(20, 5).Value = 1.037254103397646
$ws.Cells.Item(20, 6).Value = 1.050916118366013
$ws.Cells.Item(20, 9).Value = 1.041951574363608
$ws.Cells.Item(20, 10).Value = 1.039653727767642
$ws.Cells.Item(20, 11).Value = 1.045674976311071
$ws.Cells.Item(20, 12).Value = 1.040623268531943
$ws.Cells.Item(20, 13).Value = 1.054238362720516
$ws.Cells.Item(20, 14).Value = 1.01716342247978

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.032506598984268
$ws.Cells.Item(21, 4).Value = 1.041595380406331
$ws.Cells.Item(21, 5).Value = 1.036382030942346
$ws.Cells.Item(21, 6).Value = 1.050035060260693
$ws.Cells.Item(21, 9).Value = 1.041690258065641
$ws.Cells.Item(21, 10).Value = 1.039117647214078
$ws.Cells.Item(21, 11).Value = 1.0451548184579
$ws.Cells.Item(21, 12).Value = 1.03996074787949
$ws.Cells.Item(21, 13).Value = 1.053563736775113
$ws.Cells.Item(21, 14).Value = 1.016984039020935

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.03192022269129
$ws.Cells.Item(22, 4).Value = 1.041138249403111
$ws.Cells.Item(22, 5).Value = 1.035834786164171
$ws.Cells.Item(22, 6).Value = 1.049482065153228
$ws.Cells.Item(22, 9).Value = 1.041524928129584
$ws.Cells.Item(22, 10).Value = 1.038780571766743
$ws.Cells.Item(22, 11).Value = 1.044827618784125
$ws.Cells.Item(22, 12).Value = 1.039544543646379
$ws.Cells.Item(22, 13).Value = 1.053139820536489
$ws.Cells.Item(22, 14).Value = 1.016871211553483

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.032231004677034
$ws.Cells.Item(23, 4).Value = 1.041380526603251
$ws.Cells.Item(23, 5).Value = 1.036124789688914
$ws.Cells.Item(23, 6).Value = 1.049775126551082
$ws.Cells.Item(23, 9).Value = 1.04161267001247
$ws.Cells.Item(23, 10).Value = 1.03895926333746
$ws.Cells.Item(23, 11).Value = 1.045001087842848
$ws.Cells.Item(23, 12).Value = 1.039765147811104
$ws.Cells.Item(23, 13).Value = 1.053364522514673
$ws.Cells.Item(23, 14).Value = 1.016931027349

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.033456004630299
$ws.Cells.Item(24, 4).Value = 1.042335593858519
$ws.Cells.Item(24, 5).Value = 1.037268741665623
$ws.Cells.Item(24, 6).Value = 1.050930905565974
$ws.Cells.Item(24, 9).Value = 1.041955937843946
$ws.Cells.Item(24, 10).Value = 1.039662714789908
$ws.Cells.Item(24, 11).Value = 1.045683694086882
$ws.Cells.Item(24, 12).Value = 1.04063438157805
$ws.Cells.Item(24, 13).Value = 1.054249676995133
$ws.Cells.Item(24, 14).Value = 1.017166429120906

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.034880553023315
$ws.Cells.Item(25, 4).Value = 1.043446406290431
$ws.Cells.Item(25, 5).Value = 1.038600793081284
$ws.Cells.Item(25, 6).Value = 1.052276240126037
$ws.Cells.Item(25, 9).Value = 1.042349768207065
$ws.Cells.Item(25, 10).Value = 1.040478899440497
$ws.Cells.Item(25, 11).Value = 1.046475097411886
$ws.Cells.Item(25, 12).Value = 1.041644546205295
$ws.Cells.Item(25, 13).Value = 1.05527787389293
$ws.Cells.Item(25, 14).Value = 1.017439401151539
